$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.098.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5221"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2601"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06311"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.69"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07682"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.634.45"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.412"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.863.27"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5551"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8188"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.06"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.096.89"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.712"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.03"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.174"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.43"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.412"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1205"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.82"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.387"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05900"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.258"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.442"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.404"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.649"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9842"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.754"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5653"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.11%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8531"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.700"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.026.69"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.16"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.790.78"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.75"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.090"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05149"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.73%  "
